$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price cells so Excel keeps them as text
$textCells = @("D5", "D6", "D10", "D13", "D15", "D20", "D21", "D23", "D24", "D29", "D32", "D35", "D36", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "59.459.65"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "2.525.11"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "536.57"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "139.85"
$ws.Range("E6").Value = "  -3.20%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").Value = "2.529.63"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").Value = "0.0995"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").Value = "2.972.62"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "23.16"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").Value = "59.393.49"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "2.510.70"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("D20").Value = "4.23"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").Value = "321.97"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "5.81"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").Value = "61.40"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("E25").Value = "  -4.12%  "
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").Value = "6.74"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").Value = "0.0₃0769"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").Value = "160.62"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").Value = "1.12"
$ws.Range("E35").Value = "  -7.45%  "
$ws.Range("D36").Value = "18.56"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  -4.33%  "
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("D39").Value = "37.00"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").Value = "3.66"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("D41").Value = "0.808"
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "285.35"
$ws.Range("E42").Value = "  -6.12%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "5.26"
$ws.Range("E43").Value = "  -7.75%  "
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "10.87"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.599"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").Value = "124.05"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("D48").Value = "0.0925"
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("D49").Value = "18.58"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("E51").Value = "  -2.05%  "

# Restore default (Normal) style so no stray formatting is introduced
foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}
